$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")
$ws.Activate()
Write-Host $ws.Range("B26").Value
Write-Host $ws.Range("B27").Value
